$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values to reduced (custom) accuracy ---
$ws.Range("B5").Value  = 10.57
$ws.Range("C5").Value  = 7.6
$ws.Range("D5").Value  = 0.91
$ws.Range("E5").Value  = 22.68
$ws.Range("F5").Value  = 18.77
$ws.Range("G5").Value  = 8.32
$ws.Range("H5").Value  = 31.22
$ws.Range("I5").Value  = 12.8
$ws.Range("J5").Value  = 5.58
$ws.Range("K5").Value  = 8.4
$ws.Range("L5").Value  = 9.18
$ws.Range("M5").Value  = 9.54
$ws.Range("N5").Value  = 2.65
$ws.Range("O5").Value  = 8.27
$ws.Range("P5").Value  = 11.67
$ws.Range("Q5").Value  = 7.08
$ws.Range("R5").Value  = 0.78
$ws.Range("S5").Value  = 0.53
$ws.Range("T5").Value  = 118.55
$ws.Range("U5").Value  = 23.04
$ws.Range("V5").Value  = 7.63
$ws.Range("W5").Value  = 15.34
$ws.Range("X5").Value  = 8.28
$ws.Range("Y5").Value  = 0.99
$ws.Range("Z5").Value  = 14.94
$ws.Range("AA5").Value = 6.74
$ws.Range("AB5").Value = 6.06
$ws.Range("AC5").Value = 7.11
$ws.Range("AD5").Value = 9.58
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 27.89
$ws.Range("AG5").Value = 4.25
$ws.Range("AH5").Value = 9.54

# --- Remove the last data row (row 6) entirely ---
$ws.Rows("6:6").Delete()

# --- Narrow column B (width 8 -> 7 character units) ---
# ColumnWidth uses Excel's character-width unit, which this engine stores
# internally with a constant +5/6 offset; subtract it so the persisted
# <col width="..."> lands exactly on 7.
$ws.Columns("B:B").ColumnWidth = 7 - (5/6)
